# Daily update at 8 AM UTC
# - Row 60 (previously the "latest" row) reverts to the standard
#   date/time number format used by all prior rows.
# - A new row 61 is appended with the new day's data and takes over
#   the "latest row" date-only number format that row 60 used to have.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60 goes back to the regular timestamp format shared by rows 2-59.
$ws.Cells.Item(60, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data in row 61.
$ws.Cells.Item(61, 1).Value = 45647
$ws.Cells.Item(61, 2).Value = 143
$ws.Cells.Item(61, 3).Value = 133
$ws.Cells.Item(61, 4).Value = 141

# Row 61 becomes the new "latest" row, using the date-only format.
$ws.Cells.Item(61, 1).NumberFormat = "YYYY-MM-DD"
